$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the relevant paragraphs by their text, walking the document
# structurally (rather than hard-coding paragraph indices) so the
# script is resilient to anything earlier in the document.
#
#   ... "Add Animations"   (ilvl 2)
#         "Running"        (ilvl 3)  <- unchanged, already highlighted
#         "Jumping"        (ilvl 3)  <- newly added Jump animation task
#         "Falling"        (ilvl 3)  <- newly added Fall animation task
#   ... "Add Appropriate sound effects to the cat when" (ilvl 2)
# ------------------------------------------------------------------

$addAnimationsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd("`r") -eq "Add Animations") {
        $addAnimationsPara = $candidate
        break
    }
}

$runningPara = $addAnimationsPara.Next()
$jumpingPara = $runningPara.Next()
$fallingPara = $jumpingPara.Next()

# Sanity-check we found the right paragraphs before mutating anything.
if ($jumpingPara.Range.Text.TrimEnd("`r") -ne "Jumping") {
    throw "Expected 'Jumping' paragraph, found '$($jumpingPara.Range.Text)'"
}
if ($fallingPara.Range.Text.TrimEnd("`r") -ne "Falling") {
    throw "Expected 'Falling' paragraph, found '$($fallingPara.Range.Text)'"
}

# ------------------------------------------------------------------
# Highlight the new Jump-animation sub-tasks in cyan, same as the
# other already-highlighted siblings ("Running", "Add Appropriate
# sound effects..."). Using Range.Font (rather than Range directly)
# also paints the paragraph mark's run properties, matching how the
# existing cyan-highlighted paragraphs in this document are stored.
# ------------------------------------------------------------------

$addAnimationsPara.Range.Font.HighlightColorIndex = "cyan"
$jumpingPara.Range.Font.HighlightColorIndex = "cyan"
$fallingPara.Range.Font.HighlightColorIndex = "cyan"

# ------------------------------------------------------------------
# Move the "_GoBack" bookmark (Word's "last edit location" marker)
# from the end of the "Add Appropriate sound effects..." paragraph
# to the end of the "Falling" paragraph, right after its text.
# ------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-resolve the "Falling" paragraph's text end (excluding the
# trailing paragraph mark).
$fallingTextRange = $fallingPara.Range.Duplicate()
$fallingTextRange.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1
$insertionPoint = $fallingTextRange.End

# Placing a bookmark exactly at a paragraph's text/mark boundary is
# unreliable, so briefly insert a placeholder character right there,
# add the bookmark against that now-unambiguous position, then
# remove the placeholder again. The bookmark stays correctly anchored
# immediately after "Falling".
$placeholder = $d.Range($insertionPoint, $insertionPoint)
$placeholder.InsertAfter("X")

$bookmarkRange = $d.Range($insertionPoint, $insertionPoint)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$d.Range($insertionPoint, $insertionPoint + 1).Delete()
